$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 is the "古村落" (village) dungeon. Update its stat spread and add a
# new quest entry ("basement;1") to its QuestDungeon list.
$ws.Range("G6").Value = 2
$ws.Range("H6").Value = 3
$ws.Range("I6").Value = 2
$ws.Range("L6").Value = "trees;2|sandland;2|potteryroom;2|honeyhome;2|snare;1|basement;1"

# Move the active selection to H6.
$ws.Range("H6").Select()
